# Auto update stock data
# Applies refreshed "as of" row data across all five ticker sheets:
#   - A2 ("Current ..." label) is cleared
#   - A3:A8 period labels are collapsed to their compact numeric form
#   - B2 (current EBITDA figure) is refreshed for Reliance Steel & Kaiser Aluminum
#   - G2:G8 (Altman Z-Score) is filled in for Kaiser Aluminum

$wb = $excel.ActiveWorkbook

# Compact date-label replacements shared by every sheet (rows 3-8 of column A).
$dateLabels = @{
    3 = "202424312024"
    4 = "202323312023"
    5 = "202222312022"
    6 = "202121312021"
    7 = "202020312020"
    8 = "2015201920152019"
}

foreach ($ws in $wb.Worksheets) {
    # Clear the "Current ..." label in A2.
    $ws.Range("A2").Value = ""

    # Force column A (text labels) to stay text, then rewrite the period labels.
    foreach ($row in $dateLabels.Keys) {
        $cell = $ws.Cells.Item($row, 1)
        $cell.Value = "'" + $dateLabels[$row]
    }
}

# Reliance Steel & Aluminum: refresh current EBITDA figure (kept as text).
$wsRS = $wb.Worksheets.Item("Reliance Steel & Aluminum")
$wsRS.Range("B2").Value = "'12.11"

# Kaiser Aluminum: refresh current EBITDA figure and populate Altman Z-Score column.
$wsKALU = $wb.Worksheets.Item("Kaiser Aluminum")
$wsKALU.Range("B2").Value = "'9.94"
for ($row = 2; $row -le 8; $row++) {
    $wsKALU.Cells.Item($row, 7).Value = 2.22
}

Write-Host "Stock data refreshed"
